$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("main")
$v = $ws.Cells.Item(156, 2).Value()
Write-Host $v
$v2 = $ws.Range("C156").Value()
Write-Host $v2
$ws.Range("B156").Value = 1
$v3 = $ws.Range("B156").Value()
Write-Host $v3
